$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.770.93'
$ws.Range('E2').Value = '  -2.33%  '
$ws.Range('D3').Value = '1.751.41'
$ws.Range('E3').Value = '  -4.40%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = "'236.95"
$ws.Range('E5').Value = '  -6.43%  '
$ws.Range('D6').Value = "'1.000"
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('D7').Value = "'0.5073"
$ws.Range('E7').Value = '  -3.72%  '
$ws.Range('D8').Value = "'41.54"
$ws.Range('E8').Value = '  -6.54%  '
$ws.Range('D9').Value = "'0.2643"
$ws.Range('E9').Value = '  -5.23%  '
$ws.Range('D10').Value = "'0.06155"
$ws.Range('E10').Value = '  -10.47%  '
$ws.Range('D11').Value = '1.752.42'
$ws.Range('E11').Value = '  -4.54%  '
$ws.Range('D12').Value = "'15.78"
$ws.Range('E12').Value = '  -4.38%  '
$ws.Range('D13').Value = "'0.06909"
$ws.Range('E13').Value = '  -2.90%  '
$ws.Range('D14').Value = "'0.6070"
$ws.Range('E14').Value = '  -11.19%  '
$ws.Range('D15').Value = "'4.506"
$ws.Range('E15').Value = '  -7.21%  '
$ws.Range('D16').Value = "'77.16"
$ws.Range('E16').Value = '  -10.82%  '
$ws.Range('E17').Value = '  -0.25%  '
$ws.Range('D18').Value = "'1.000"
$ws.Range('E18').Value = '  -0.16%  '
$ws.Range('D19').Value = '25.781.90'
$ws.Range('E19').Value = '  -2.34%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = "'0.000006867"
$ws.Range('E20').Value = '  -6.70%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').Value = "'11.72"
$ws.Range('E21').Value = '  -11.16%  '
$ws.Range('D22').Value = '1.973.66'
$ws.Range('E22').Value = '  -5.56%  '
$ws.Range('D23').Value = "'4.094"
$ws.Range('E23').Value = '  -8.92%  '
$ws.Range('D24').Value = "'8.261"
$ws.Range('E24').Value = '  -7.82%  '
$ws.Range('D25').Value = "'5.218"
$ws.Range('E25').Value = '  -9.97%  '
$ws.Range('D26').Value = "'137.65"
$ws.Range('E26').Value = '  -3.13%  '
$ws.Range('E27').Value = '  -12.25%  '
$ws.Range('D28').Value = "'1.831"
$ws.Range('E28').Value = '  -9.58%  '
$ws.Range('D29').Value = "'15.04"
$ws.Range('E29').Value = '  -9.08%  '
$ws.Range('D30').Value = "'102.85"
$ws.Range('E30').Value = '  -5.64%  '
$ws.Range('D31').Value = "'0.08211"
$ws.Range('E31').Value = '  -5.84%  '
$ws.Range('D32').Value = "'3.698"
$ws.Range('E32').Value = '  -9.16%  '
$ws.Range('D33').Value = "'3.467"
$ws.Range('E33').Value = '  -9.97%  '
$ws.Range('D34').Value = "'0.04516"
$ws.Range('E34').Value = '  -3.99%  '
$ws.Range('D35').Value = "'0.9998"
$ws.Range('E35').Value = '  -0.12%  '
$ws.Range('D36').Value = "'2.666"
$ws.Range('E36').Value = '  -7.69%  '
$ws.Range('D37').Value = "'0.9998"
$ws.Range('E37').Value = '  -9.87%  '
$ws.Range('D38').Value = "'0.6093"
$ws.Range('E38').Value = '  -13.35%  '
$ws.Range('D39').Value = "'2.693"
$ws.Range('E39').Value = '  -11.97%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = "'1.955"
$ws.Range('E40').Value = '  -9.98%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = "'0.01555"
$ws.Range('E41').Value = '  -4.75%  '
$ws.Range('D42').Value = "'1.000"
$ws.Range('E42').Value = '  -0.14%  '
$ws.Range('D43').Value = "'103.65"
$ws.Range('E43').Value = '  -1.37%  '
$ws.Range('D44').Value = "'0.3834"
$ws.Range('E44').Value = '  -13.55%  '
$ws.Range('D45').Value = "'0.7396"
$ws.Range('E45').Value = '  -14.23%  '
$ws.Range('D46').Value = "'4.930"
$ws.Range('E46').Value = '  -13.80%  '
$ws.Range('D47').Value = "'0.05468"
$ws.Range('E47').Value = '  -1.97%  '
$ws.Range('D48').Value = "'0.1104"
$ws.Range('E48').Value = '  -6.68%  '
$ws.Range('D49').Value = "'6.012"
$ws.Range('E49').Value = '  -14.07%  '
$ws.Range('D50').Value = "'7.701"
$ws.Range('E50').Value = '  -10.82%  '
$ws.Range('D51').Value = "'29.98"
$ws.Range('E51').Value = '  -10.17%  '
